$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("A24").Value = "IPA113"
$ws.Range("B24").Value = "OBT"
$ws.Range("C24").Value = "Save the technology search data and rerun the saved data"
$ws.Range("D24").Value = "Y"
